$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D, E, G (Price, Volume(1h), Hora) to Text format first so that
# numeric-looking / percent-looking strings round-trip as literal text, matching
# the source workbook (which stores these as inline text, not numbers).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Apply the updated scraped values (price / 1h volume % / hour) per the commit diff.
$ws.Range("D2").Value = "277.53"
$ws.Range("E2").Value = "1.59%"
$ws.Range("G2").Value = "21"
$ws.Range("D3").Value = "27.32"
$ws.Range("E3").Value = "1.93%"
$ws.Range("G3").Value = "21"
$ws.Range("D4").Value = "4.870"
$ws.Range("E4").Value = "-0.82%"
$ws.Range("G4").Value = "21"
$ws.Range("D5").Value = "0.06366"
$ws.Range("G5").Value = "21"
$ws.Range("D6").Value = "6.973"
$ws.Range("E6").Value = "0.55%"
$ws.Range("G6").Value = "21"
$ws.Range("D7").Value = "1.242"
$ws.Range("E7").Value = "-7.38%"
$ws.Range("G7").Value = "21"
$ws.Range("D8").Value = "0.8821"
$ws.Range("E8").Value = "-0.38%"
$ws.Range("G8").Value = "21"
$ws.Range("D9").Value = "0.1528"
$ws.Range("E9").Value = "3.79%"
$ws.Range("G9").Value = "21"
$ws.Range("D10").Value = "0.05106"
$ws.Range("E10").Value = "-0.45%"
$ws.Range("G10").Value = "21"
$ws.Range("D11").Value = "0.07566"
$ws.Range("E11").Value = "2.51%"
$ws.Range("G11").Value = "21"
$ws.Range("D12").Value = "0.02979"
$ws.Range("E12").Value = "-5.85%"
$ws.Range("G12").Value = "21"
$ws.Range("D13").Value = "0.09010"
$ws.Range("E13").Value = "-0.61%"
$ws.Range("G13").Value = "21"
$ws.Range("D14").Value = "0.001568"
$ws.Range("E14").Value = "-0.77%"
$ws.Range("G14").Value = "21"
$ws.Range("D15").Value = "0.0006419"
$ws.Range("E15").Value = "1.50%"
$ws.Range("G15").Value = "21"
$ws.Range("D16").Value = "0.005918"
$ws.Range("E16").Value = "-2.07%"
$ws.Range("G16").Value = "21"
$ws.Range("E17").Value = "-0.48%"
$ws.Range("G17").Value = "21"
$ws.Range("D18").Value = "3.313"
$ws.Range("E18").Value = "-1.24%"
$ws.Range("G18").Value = "21"
$ws.Range("E19").Value = "-0.35%"
$ws.Range("G19").Value = "21"
$ws.Range("G20").Value = "21"
$ws.Range("D21").Value = "0.1336"
$ws.Range("E21").Value = "0.21%"
$ws.Range("G21").Value = "21"
$ws.Range("D22").Value = "3.899"
$ws.Range("E22").Value = "-0.13%"
$ws.Range("G22").Value = "21"
$ws.Range("D23").Value = "0.04418"
$ws.Range("E23").Value = "1.67%"
$ws.Range("G23").Value = "21"
$ws.Range("D24").Value = "0.001170"
$ws.Range("E24").Value = "-0.89%"
$ws.Range("G24").Value = "21"
$ws.Range("D25").Value = "0.003873"
$ws.Range("E25").Value = "6.27%"
$ws.Range("G25").Value = "21"
$ws.Range("E26").Value = "-0.25%"
$ws.Range("G26").Value = "21"
$ws.Range("E27").Value = "-0.31%"
$ws.Range("G27").Value = "21"
$ws.Range("G28").Value = "21"
$ws.Range("G29").Value = "21"
$ws.Range("G30").Value = "21"
$ws.Range("G31").Value = "21"
$ws.Range("G32").Value = "21"
$ws.Range("G33").Value = "21"
$ws.Range("G34").Value = "21"
$ws.Range("G35").Value = "21"
$ws.Range("G36").Value = "21"
$ws.Range("G37").Value = "21"
$ws.Range("G38").Value = "21"
$ws.Range("G39").Value = "21"
$ws.Range("D40").Value = "0.04141"
$ws.Range("E40").Value = "2.85%"
$ws.Range("G40").Value = "21"
$ws.Range("D41").Value = "0.006826"
$ws.Range("E41").Value = "3.05%"
$ws.Range("G41").Value = "21"
$ws.Range("D42").Value = "0.1181"
$ws.Range("E42").Value = "1.30%"
$ws.Range("G42").Value = "21"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").Value = "-8.82%"
$ws.Range("G43").Value = "21"
$ws.Range("D44").Value = "0.01175"
$ws.Range("E44").Value = "-6.81%"
$ws.Range("G44").Value = "21"
$ws.Range("D45").Value = "0.00005183"
$ws.Range("E45").Value = "-1.51%"
$ws.Range("G45").Value = "21"
$ws.Range("G46").Value = "21"
$ws.Range("E47").Value = "-4.83%"
$ws.Range("G47").Value = "21"
$ws.Range("G48").Value = "21"
$ws.Range("G49").Value = "21"
$ws.Range("G50").Value = "21"
$ws.Range("G51").Value = "21"
